$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C8").Value = "Ehab"
$ws.Range("C8").Select()
